$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.1
$ws.Range("K2").Value = 2.3
$ws.Range("N2").Value = 12
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.36
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("X2").Value = 7.5
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 13
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 17
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.5
$ws.Range("AQ2").Value = 23
$ws.Range("AT2").Value = 3
$ws.Range("AX2").Value = 29
$ws.Range("AY2").Value = 34
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 126
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 1.8
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.5
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("W4").Value = 12
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 17
$ws.Range("Z4").Value = 51
$ws.Range("AH4").Value = 6
$ws.Range("AI4").Value = 7.5
$ws.Range("AK4").Value = 13
$ws.Range("AL4").Value = 15
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 6.5
$ws.Range("AO4").Value = 26
$ws.Range("AQ4").Value = 101
$ws.Range("AT4").Value = 2.63
$ws.Range("AU4").Value = 9
$ws.Range("AW4").Value = 3.6
$ws.Range("AX4").Value = 9.5
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 4.33
$ws.Range("I5").Value = 1.48
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 2.38
$ws.Range("L5").Value = 2.05
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1.75
$ws.Range("R5").Value = 2.05
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 3.25
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.83
$ws.Range("W5").Value = 15
$ws.Range("X5").Value = 29
$ws.Range("Y5").Value = 17
$ws.Range("Z5").Value = 67
$ws.Range("AA5").Value = 41
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 13
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 19
$ws.Range("AG5").Value = 301
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 7.5
$ws.Range("AK5").Value = 11
$ws.Range("AL5").Value = 12
$ws.Range("AN5").Value = 7.5
$ws.Range("AO5").Value = 34
$ws.Range("AP5").Value = 34
$ws.Range("AQ5").Value = 126
$ws.Range("AR5").Value = 126
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 3.25
$ws.Range("AU5").Value = 8.5
$ws.Range("AW5").Value = 3.5
$ws.Range("AX5").Value = 7.5
$ws.Range("AZ5").Value = 21
$ws.Range("BA5").Value = 41
$ws.Range("N6").Value = 9
